$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 208 <- old row 212 (id 7404214)
$ws.Cells.Item(208, 2).Value = 7404214
$ws.Cells.Item(208, 5).Value = "Boyaca Chico"
$ws.Cells.Item(208, 6).Value = "Deportivo Cali"
$ws.Cells.Item(208, 7).Value = 1
$ws.Cells.Item(208, 8).Value = 1
$ws.Cells.Item(208, 9).Value = 0
$ws.Cells.Item(208, 10).Value = 0
$ws.Cells.Item(208, 11).Value = "D"
$ws.Cells.Item(208, 12).Value = 3.2
$ws.Cells.Item(208, 13).Value = 3.1
$ws.Cells.Item(208, 14).Value = 2.2
$ws.Cells.Item(208, 15).Value = 3.6
$ws.Cells.Item(208, 16).Value = 3
$ws.Cells.Item(208, 17).Value = 2.25
$ws.Cells.Item(208, 18).Value = 0.25
$ws.Cells.Item(208, 19).Value = 1.95
$ws.Cells.Item(208, 20).Value = 1.9
$ws.Cells.Item(208, 21).Value = 2.25
$ws.Cells.Item(208, 22).Value = 1.875
$ws.Cells.Item(208, 23).Value = 1.975
$ws.Cells.Item(208, 24).Value = -1
$ws.Cells.Item(208, 25).Value = 2
$ws.Cells.Item(208, 26).Value = -1
$ws.Cells.Item(208, 27).Value = 0.475
$ws.Cells.Item(208, 28).Value = -0.5
$ws.Cells.Item(208, 29).Value = -0.5
$ws.Cells.Item(208, 30).Value = 0.4875

# Row 209 <- old row 208 (id 7404217)
$ws.Cells.Item(209, 2).Value = 7404217
$ws.Cells.Item(209, 5).Value = "Alianza Petrolera"
$ws.Cells.Item(209, 6).Value = "Deportivo Pereira"
$ws.Cells.Item(209, 7).Value = 2
$ws.Cells.Item(209, 8).Value = 1
$ws.Cells.Item(209, 9).Value = 2
$ws.Cells.Item(209, 10).Value = 1
$ws.Cells.Item(209, 11).Value = "H"
$ws.Cells.Item(209, 12).Value = 1.95
$ws.Cells.Item(209, 13).Value = 3.2
$ws.Cells.Item(209, 14).Value = 3.75
$ws.Cells.Item(209, 15).Value = 1.95
$ws.Cells.Item(209, 16).Value = 3.2
$ws.Cells.Item(209, 17).Value = 4.75
$ws.Cells.Item(209, 18).Value = -0.5
$ws.Cells.Item(209, 19).Value = 1.925
$ws.Cells.Item(209, 20).Value = 1.875
$ws.Cells.Item(209, 21).Value = 2
$ws.Cells.Item(209, 22).Value = 1.825
$ws.Cells.Item(209, 23).Value = 1.975
$ws.Cells.Item(209, 24).Value = 0.95
$ws.Cells.Item(209, 25).Value = -1
$ws.Cells.Item(209, 26).Value = -1
$ws.Cells.Item(209, 27).Value = 0.925
$ws.Cells.Item(209, 28).Value = -1
$ws.Cells.Item(209, 29).Value = 0.825
$ws.Cells.Item(209, 30).Value = -1

# Row 210 <- old row 209 (id 7404218)
$ws.Cells.Item(210, 2).Value = 7404218
$ws.Cells.Item(210, 5).Value = "Junior"
$ws.Cells.Item(210, 6).Value = "Atletico Huila"
$ws.Cells.Item(210, 7).Value = 2
$ws.Cells.Item(210, 8).Value = 0
$ws.Cells.Item(210, 9).Value = 1
$ws.Cells.Item(210, 10).Value = 0
$ws.Cells.Item(210, 11).Value = "H"
$ws.Cells.Item(210, 12).Value = 1.363
$ws.Cells.Item(210, 13).Value = 4.5
$ws.Cells.Item(210, 14).Value = 7
$ws.Cells.Item(210, 15).Value = 1.3
$ws.Cells.Item(210, 16).Value = 5
$ws.Cells.Item(210, 17).Value = 12
$ws.Cells.Item(210, 18).Value = -1.5
$ws.Cells.Item(210, 19).Value = 1.9
$ws.Cells.Item(210, 20).Value = 1.95
$ws.Cells.Item(210, 21).Value = 2.75
$ws.Cells.Item(210, 22).Value = 2.025
$ws.Cells.Item(210, 23).Value = 1.825
$ws.Cells.Item(210, 24).Value = 0.3
$ws.Cells.Item(210, 25).Value = -1
$ws.Cells.Item(210, 26).Value = -1
$ws.Cells.Item(210, 27).Value = 0.8999999999999999
$ws.Cells.Item(210, 28).Value = -1
$ws.Cells.Item(210, 29).Value = -1
$ws.Cells.Item(210, 30).Value = 0.825

# Row 211 <- old row 210 (id 7404212)
$ws.Cells.Item(211, 2).Value = 7404212
$ws.Cells.Item(211, 5).Value = "Envigado FC"
$ws.Cells.Item(211, 6).Value = "Deportivo Pasto"
$ws.Cells.Item(211, 7).Value = 1
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(211, 9).Value = 1
$ws.Cells.Item(211, 10).Value = 1
$ws.Cells.Item(211, 11).Value = "D"
$ws.Cells.Item(211, 12).Value = 2.6
$ws.Cells.Item(211, 13).Value = 2.875
$ws.Cells.Item(211, 14).Value = 2.8
$ws.Cells.Item(211, 15).Value = 2.8
$ws.Cells.Item(211, 16).Value = 3.2
$ws.Cells.Item(211, 17).Value = 2.625
$ws.Cells.Item(211, 18).Value = 0
$ws.Cells.Item(211, 19).Value = 1.975
$ws.Cells.Item(211, 20).Value = 1.875
$ws.Cells.Item(211, 21).Value = 2.5
$ws.Cells.Item(211, 22).Value = 2.025
$ws.Cells.Item(211, 23).Value = 1.825
$ws.Cells.Item(211, 24).Value = -1
$ws.Cells.Item(211, 25).Value = 2.2
$ws.Cells.Item(211, 26).Value = -1
$ws.Cells.Item(211, 27).Value = 0
$ws.Cells.Item(211, 28).Value = 0
$ws.Cells.Item(211, 29).Value = -1
$ws.Cells.Item(211, 30).Value = 0.825

# Row 212 <- old row 211 (id 7404216)
$ws.Cells.Item(212, 2).Value = 7404216
$ws.Cells.Item(212, 5).Value = "Independiente Santa Fe"
$ws.Cells.Item(212, 6).Value = "Once Caldas"
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 1
$ws.Cells.Item(212, 9).Value = 0
$ws.Cells.Item(212, 10).Value = 1
$ws.Cells.Item(212, 11).Value = "A"
$ws.Cells.Item(212, 12).Value = 1.85
$ws.Cells.Item(212, 13).Value = 3.1
$ws.Cells.Item(212, 14).Value = 4.2
$ws.Cells.Item(212, 15).Value = 2.25
$ws.Cells.Item(212, 16).Value = 3.2
$ws.Cells.Item(212, 17).Value = 3.3
$ws.Cells.Item(212, 18).Value = -0.25
$ws.Cells.Item(212, 19).Value = 1.9
$ws.Cells.Item(212, 20).Value = 1.9
$ws.Cells.Item(212, 21).Value = 2.5
$ws.Cells.Item(212, 22).Value = 1.925
$ws.Cells.Item(212, 23).Value = 1.925
$ws.Cells.Item(212, 24).Value = -1
$ws.Cells.Item(212, 25).Value = -1
$ws.Cells.Item(212, 26).Value = 2.3
$ws.Cells.Item(212, 27).Value = -1
$ws.Cells.Item(212, 28).Value = 0.8999999999999999
$ws.Cells.Item(212, 29).Value = -1
$ws.Cells.Item(212, 30).Value = 0.925

# Row 213 <- old row 216 (id 7404260)
$ws.Cells.Item(213, 2).Value = 7404260
$ws.Cells.Item(213, 5).Value = "Atletico Nacional Medellin"
$ws.Cells.Item(213, 6).Value = "Deportes Tolima"
$ws.Cells.Item(213, 7).Value = 2
$ws.Cells.Item(213, 8).Value = 3
$ws.Cells.Item(213, 9).Value = 0
$ws.Cells.Item(213, 10).Value = 1
$ws.Cells.Item(213, 11).Value = "A"
$ws.Cells.Item(213, 12).Value = 2
$ws.Cells.Item(213, 13).Value = 3.25
$ws.Cells.Item(213, 14).Value = 3.5
$ws.Cells.Item(213, 15).Value = 1.75
$ws.Cells.Item(213, 16).Value = 3.6
$ws.Cells.Item(213, 17).Value = 4.75
$ws.Cells.Item(213, 18).Value = -0.75
$ws.Cells.Item(213, 19).Value = 2
$ws.Cells.Item(213, 20).Value = 1.8
$ws.Cells.Item(213, 21).Value = 2.5
$ws.Cells.Item(213, 22).Value = 2
$ws.Cells.Item(213, 23).Value = 1.8
$ws.Cells.Item(213, 24).Value = -1
$ws.Cells.Item(213, 25).Value = -1
$ws.Cells.Item(213, 26).Value = 3.75
$ws.Cells.Item(213, 27).Value = -1
$ws.Cells.Item(213, 28).Value = 0.8
$ws.Cells.Item(213, 29).Value = 1
$ws.Cells.Item(213, 30).Value = -1

# Row 214 <- old row 217 (id 7404522)
$ws.Cells.Item(214, 2).Value = 7404522
$ws.Cells.Item(214, 5).Value = "La Equidad"
$ws.Cells.Item(214, 6).Value = "Millonarios"
$ws.Cells.Item(214, 7).Value = 2
$ws.Cells.Item(214, 8).Value = 1
$ws.Cells.Item(214, 9).Value = 0
$ws.Cells.Item(214, 10).Value = 1
$ws.Cells.Item(214, 11).Value = "H"
$ws.Cells.Item(214, 12).Value = 2.4
$ws.Cells.Item(214, 13).Value = 3.1
$ws.Cells.Item(214, 14).Value = 2.875
$ws.Cells.Item(214, 15).Value = 2.1
$ws.Cells.Item(214, 16).Value = 3.1
$ws.Cells.Item(214, 17).Value = 3.8
$ws.Cells.Item(214, 18).Value = -0.25
$ws.Cells.Item(214, 19).Value = 1.75
$ws.Cells.Item(214, 20).Value = 2.05
$ws.Cells.Item(214, 21).Value = 2
$ws.Cells.Item(214, 22).Value = 1.85
$ws.Cells.Item(214, 23).Value = 1.95
$ws.Cells.Item(214, 24).Value = 1.1
$ws.Cells.Item(214, 25).Value = -1
$ws.Cells.Item(214, 26).Value = -1
$ws.Cells.Item(214, 27).Value = 0.75
$ws.Cells.Item(214, 28).Value = -1
$ws.Cells.Item(214, 29).Value = 0.8500000000000001
$ws.Cells.Item(214, 30).Value = -1

# Row 216 <- old row 214 (id 7404215)
$ws.Cells.Item(216, 2).Value = 7404215
$ws.Cells.Item(216, 5).Value = "America de Cali"
$ws.Cells.Item(216, 6).Value = "Atletico Bucaramanga"
$ws.Cells.Item(216, 7).Value = 1
$ws.Cells.Item(216, 8).Value = 2
$ws.Cells.Item(216, 9).Value = 1
$ws.Cells.Item(216, 10).Value = 1
$ws.Cells.Item(216, 11).Value = "A"
$ws.Cells.Item(216, 12).Value = 1.444
$ws.Cells.Item(216, 13).Value = 4.5
$ws.Cells.Item(216, 14).Value = 6
$ws.Cells.Item(216, 15).Value = 1.363
$ws.Cells.Item(216, 16).Value = 5
$ws.Cells.Item(216, 17).Value = 7.5
$ws.Cells.Item(216, 18).Value = -1.25
$ws.Cells.Item(216, 19).Value = 1.775
$ws.Cells.Item(216, 20).Value = 2.025
$ws.Cells.Item(216, 21).Value = 3
$ws.Cells.Item(216, 22).Value = 1.925
$ws.Cells.Item(216, 23).Value = 1.875
$ws.Cells.Item(216, 24).Value = -1
$ws.Cells.Item(216, 25).Value = -1
$ws.Cells.Item(216, 26).Value = 6.5
$ws.Cells.Item(216, 27).Value = -1
$ws.Cells.Item(216, 28).Value = 1.025
$ws.Cells.Item(216, 29).Value = 0
$ws.Cells.Item(216, 30).Value = 0

# Row 217 <- old row 213 (id 7404213)
$ws.Cells.Item(217, 2).Value = 7404213
$ws.Cells.Item(217, 5).Value = "Jaguares de Cordoba"
$ws.Cells.Item(217, 6).Value = "Aguilas Doradas"
$ws.Cells.Item(217, 7).Value = 0
$ws.Cells.Item(217, 8).Value = 1
$ws.Cells.Item(217, 9).Value = 0
$ws.Cells.Item(217, 10).Value = 0
$ws.Cells.Item(217, 11).Value = "A"
$ws.Cells.Item(217, 12).Value = 3.25
$ws.Cells.Item(217, 13).Value = 3.1
$ws.Cells.Item(217, 14).Value = 2.2
$ws.Cells.Item(217, 15).Value = 3.6
$ws.Cells.Item(217, 16).Value = 3.2
$ws.Cells.Item(217, 17).Value = 2.15
$ws.Cells.Item(217, 18).Value = 0.25
$ws.Cells.Item(217, 19).Value = 1.975
$ws.Cells.Item(217, 20).Value = 1.825
$ws.Cells.Item(217, 21).Value = 2
$ws.Cells.Item(217, 22).Value = 1.75
$ws.Cells.Item(217, 23).Value = 2.05
$ws.Cells.Item(217, 24).Value = -1
$ws.Cells.Item(217, 25).Value = -1
$ws.Cells.Item(217, 26).Value = 1.15
$ws.Cells.Item(217, 27).Value = -1
$ws.Cells.Item(217, 28).Value = 0.825
$ws.Cells.Item(217, 29).Value = -1
$ws.Cells.Item(217, 30).Value = 1.05

# Row 240 <- old row 241 (id 7528135)
$ws.Cells.Item(240, 2).Value = 7528135
$ws.Cells.Item(240, 5).Value = "Independiente Medellin"
$ws.Cells.Item(240, 6).Value = "America de Cali"
$ws.Cells.Item(240, 7).Value = 2
$ws.Cells.Item(240, 8).Value = 1
$ws.Cells.Item(240, 9).Value = 2
$ws.Cells.Item(240, 10).Value = 1
$ws.Cells.Item(240, 11).Value = "H"
$ws.Cells.Item(240, 12).Value = 2.15
$ws.Cells.Item(240, 13).Value = 3.3
$ws.Cells.Item(240, 14).Value = 3.4
$ws.Cells.Item(240, 15).Value = 2.375
$ws.Cells.Item(240, 16).Value = 3.3
$ws.Cells.Item(240, 17).Value = 3.1
$ws.Cells.Item(240, 18).Value = -0.25
$ws.Cells.Item(240, 19).Value = 2
$ws.Cells.Item(240, 20).Value = 1.8
$ws.Cells.Item(240, 21).Value = 2.5
$ws.Cells.Item(240, 22).Value = 1.975
$ws.Cells.Item(240, 23).Value = 1.825
$ws.Cells.Item(240, 24).Value = 1.375
$ws.Cells.Item(240, 25).Value = -1
$ws.Cells.Item(240, 26).Value = -1
$ws.Cells.Item(240, 27).Value = 1
$ws.Cells.Item(240, 28).Value = -1
$ws.Cells.Item(240, 29).Value = 0.9750000000000001
$ws.Cells.Item(240, 30).Value = -1

# Row 241 <- old row 240 (id 7528603)
$ws.Cells.Item(241, 2).Value = 7528603
$ws.Cells.Item(241, 5).Value = "Junior"
$ws.Cells.Item(241, 6).Value = "Deportes Tolima"
$ws.Cells.Item(241, 7).Value = 4
$ws.Cells.Item(241, 8).Value = 2
$ws.Cells.Item(241, 9).Value = 3
$ws.Cells.Item(241, 10).Value = 2
$ws.Cells.Item(241, 11).Value = "H"
$ws.Cells.Item(241, 12).Value = 1.95
$ws.Cells.Item(241, 13).Value = 3.3
$ws.Cells.Item(241, 14).Value = 4
$ws.Cells.Item(241, 15).Value = 1.909
$ws.Cells.Item(241, 16).Value = 3.75
$ws.Cells.Item(241, 17).Value = 3.8
$ws.Cells.Item(241, 18).Value = -0.5
$ws.Cells.Item(241, 19).Value = 1.9
$ws.Cells.Item(241, 20).Value = 1.9
$ws.Cells.Item(241, 21).Value = 2.5
$ws.Cells.Item(241, 22).Value = 1.85
$ws.Cells.Item(241, 23).Value = 1.95
$ws.Cells.Item(241, 24).Value = 0.909
$ws.Cells.Item(241, 25).Value = -1
$ws.Cells.Item(241, 26).Value = -1
$ws.Cells.Item(241, 27).Value = 0.8999999999999999
$ws.Cells.Item(241, 28).Value = -1
$ws.Cells.Item(241, 29).Value = 0.8500000000000001
$ws.Cells.Item(241, 30).Value = -1

# Row 425 <- old row 426 (id 7658989)
$ws.Cells.Item(425, 2).Value = 7658989
$ws.Cells.Item(425, 5).Value = "Jaguares de Cordoba"
$ws.Cells.Item(425, 6).Value = "Independiente Santa Fe"
$ws.Cells.Item(425, 7).Value = 1
$ws.Cells.Item(425, 8).Value = 0
$ws.Cells.Item(425, 9).Value = 1
$ws.Cells.Item(425, 10).Value = 0
$ws.Cells.Item(425, 11).Value = "H"
$ws.Cells.Item(425, 12).Value = 3
$ws.Cells.Item(425, 13).Value = 3.2
$ws.Cells.Item(425, 14).Value = 2.3
$ws.Cells.Item(425, 15).Value = 3.4
$ws.Cells.Item(425, 16).Value = 3.6
$ws.Cells.Item(425, 17).Value = 2.05
$ws.Cells.Item(425, 18).Value = 0.25
$ws.Cells.Item(425, 19).Value = 2
$ws.Cells.Item(425, 20).Value = 1.8
$ws.Cells.Item(425, 21).Value = 2.5
$ws.Cells.Item(425, 22).Value = 1.8
$ws.Cells.Item(425, 23).Value = 2
$ws.Cells.Item(425, 24).Value = 2.4
$ws.Cells.Item(425, 25).Value = -1
$ws.Cells.Item(425, 26).Value = -1
$ws.Cells.Item(425, 27).Value = 1
$ws.Cells.Item(425, 28).Value = -1
$ws.Cells.Item(425, 29).Value = -1
$ws.Cells.Item(425, 30).Value = 1

# Row 426 <- old row 427 (id 7658914)
$ws.Cells.Item(426, 2).Value = 7658914
$ws.Cells.Item(426, 5).Value = "La Equidad"
$ws.Cells.Item(426, 6).Value = "Deportivo Pereira"
$ws.Cells.Item(426, 7).Value = 0
$ws.Cells.Item(426, 8).Value = 2
$ws.Cells.Item(426, 9).Value = 0
$ws.Cells.Item(426, 10).Value = 1
$ws.Cells.Item(426, 11).Value = "A"
$ws.Cells.Item(426, 12).Value = 2
$ws.Cells.Item(426, 13).Value = 3.1
$ws.Cells.Item(426, 14).Value = 3.75
$ws.Cells.Item(426, 15).Value = 2.25
$ws.Cells.Item(426, 16).Value = 3.2
$ws.Cells.Item(426, 17).Value = 3.3
$ws.Cells.Item(426, 18).Value = -0.25
$ws.Cells.Item(426, 19).Value = 1.925
$ws.Cells.Item(426, 20).Value = 1.875
$ws.Cells.Item(426, 21).Value = 2
$ws.Cells.Item(426, 22).Value = 1.825
$ws.Cells.Item(426, 23).Value = 1.975
$ws.Cells.Item(426, 24).Value = -1
$ws.Cells.Item(426, 25).Value = -1
$ws.Cells.Item(426, 26).Value = 2.3
$ws.Cells.Item(426, 27).Value = -1
$ws.Cells.Item(426, 28).Value = 0.875
$ws.Cells.Item(426, 29).Value = 0
$ws.Cells.Item(426, 30).Value = 0

# Row 427 <- old row 429 (id 7658915)
$ws.Cells.Item(427, 2).Value = 7658915
$ws.Cells.Item(427, 5).Value = "Once Caldas"
$ws.Cells.Item(427, 6).Value = "America de Cali"
$ws.Cells.Item(427, 7).Value = 0
$ws.Cells.Item(427, 8).Value = 0
$ws.Cells.Item(427, 9).Value = 0
$ws.Cells.Item(427, 10).Value = 0
$ws.Cells.Item(427, 11).Value = "D"
$ws.Cells.Item(427, 12).Value = 2.3
$ws.Cells.Item(427, 13).Value = 3
$ws.Cells.Item(427, 14).Value = 3.1
$ws.Cells.Item(427, 15).Value = 2.3
$ws.Cells.Item(427, 16).Value = 3.2
$ws.Cells.Item(427, 17).Value = 3.3
$ws.Cells.Item(427, 18).Value = -0.25
$ws.Cells.Item(427, 19).Value = 1.975
$ws.Cells.Item(427, 20).Value = 1.875
$ws.Cells.Item(427, 21).Value = 2.25
$ws.Cells.Item(427, 22).Value = 2.025
$ws.Cells.Item(427, 23).Value = 1.825
$ws.Cells.Item(427, 24).Value = -1
$ws.Cells.Item(427, 25).Value = 2.2
$ws.Cells.Item(427, 26).Value = -1
$ws.Cells.Item(427, 27).Value = -0.5
$ws.Cells.Item(427, 28).Value = 0.4375
$ws.Cells.Item(427, 29).Value = -1
$ws.Cells.Item(427, 30).Value = 0.825

# Row 428 <- old row 430 (id 7658988)
$ws.Cells.Item(428, 2).Value = 7658988
$ws.Cells.Item(428, 5).Value = "Envigado FC"
$ws.Cells.Item(428, 6).Value = "Independiente Medellin"
$ws.Cells.Item(428, 7).Value = 0
$ws.Cells.Item(428, 8).Value = 1
$ws.Cells.Item(428, 9).Value = 0
$ws.Cells.Item(428, 10).Value = 1
$ws.Cells.Item(428, 11).Value = "A"
$ws.Cells.Item(428, 12).Value = 4.2
$ws.Cells.Item(428, 13).Value = 3.4
$ws.Cells.Item(428, 14).Value = 1.8
$ws.Cells.Item(428, 15).Value = 5.25
$ws.Cells.Item(428, 16).Value = 3.6
$ws.Cells.Item(428, 17).Value = 1.7
$ws.Cells.Item(428, 18).Value = 0.75
$ws.Cells.Item(428, 19).Value = 1.925
$ws.Cells.Item(428, 20).Value = 1.875
$ws.Cells.Item(428, 21).Value = 2.25
$ws.Cells.Item(428, 22).Value = 1.775
$ws.Cells.Item(428, 23).Value = 2.025
$ws.Cells.Item(428, 24).Value = -1
$ws.Cells.Item(428, 25).Value = -1
$ws.Cells.Item(428, 26).Value = 0.7
$ws.Cells.Item(428, 27).Value = -0.5
$ws.Cells.Item(428, 28).Value = 0.4375
$ws.Cells.Item(428, 29).Value = -1
$ws.Cells.Item(428, 30).Value = 1.025

# Row 429 <- old row 431 (id 7736841)
$ws.Cells.Item(429, 2).Value = 7736841
$ws.Cells.Item(429, 5).Value = "Atletico Bucaramanga"
$ws.Cells.Item(429, 6).Value = "Alianza"
$ws.Cells.Item(429, 7).Value = 1
$ws.Cells.Item(429, 8).Value = 0
$ws.Cells.Item(429, 9).Value = 1
$ws.Cells.Item(429, 10).Value = 0
$ws.Cells.Item(429, 11).Value = "H"
$ws.Cells.Item(429, 12).Value = 1.666
$ws.Cells.Item(429, 13).Value = 3.5
$ws.Cells.Item(429, 14).Value = 5
$ws.Cells.Item(429, 15).Value = 1.65
$ws.Cells.Item(429, 16).Value = 3.75
$ws.Cells.Item(429, 17).Value = 5.75
$ws.Cells.Item(429, 18).Value = -0.75
$ws.Cells.Item(429, 19).Value = 1.8
$ws.Cells.Item(429, 20).Value = 2
$ws.Cells.Item(429, 21).Value = 2.25
$ws.Cells.Item(429, 22).Value = 1.9
$ws.Cells.Item(429, 23).Value = 1.9
$ws.Cells.Item(429, 24).Value = 0.6499999999999999
$ws.Cells.Item(429, 25).Value = -1
$ws.Cells.Item(429, 26).Value = -1
$ws.Cells.Item(429, 27).Value = 0.4
$ws.Cells.Item(429, 28).Value = -0.5
$ws.Cells.Item(429, 29).Value = -1
$ws.Cells.Item(429, 30).Value = 0.8999999999999999

# Row 430 <- old row 425 (id 7658990)
$ws.Cells.Item(430, 2).Value = 7658990
$ws.Cells.Item(430, 5).Value = "Millonarios"
$ws.Cells.Item(430, 6).Value = "Boyaca Chico"
$ws.Cells.Item(430, 7).Value = 3
$ws.Cells.Item(430, 8).Value = 0
$ws.Cells.Item(430, 9).Value = 1
$ws.Cells.Item(430, 10).Value = 0
$ws.Cells.Item(430, 11).Value = "H"
$ws.Cells.Item(430, 12).Value = 1.4
$ws.Cells.Item(430, 13).Value = 4.2
$ws.Cells.Item(430, 14).Value = 7
$ws.Cells.Item(430, 15).Value = 1.4
$ws.Cells.Item(430, 16).Value = 4.5
$ws.Cells.Item(430, 17).Value = 8.5
$ws.Cells.Item(430, 18).Value = -1.25
$ws.Cells.Item(430, 19).Value = 1.95
$ws.Cells.Item(430, 20).Value = 1.9
$ws.Cells.Item(430, 21).Value = 2.5
$ws.Cells.Item(430, 22).Value = 1.975
$ws.Cells.Item(430, 23).Value = 1.875
$ws.Cells.Item(430, 24).Value = 0.3999999999999999
$ws.Cells.Item(430, 25).Value = -1
$ws.Cells.Item(430, 26).Value = -1
$ws.Cells.Item(430, 27).Value = 0.95
$ws.Cells.Item(430, 28).Value = -1
$ws.Cells.Item(430, 29).Value = 0.9750000000000001
$ws.Cells.Item(430, 30).Value = -1

# Row 431 <- old row 428 (id 7658985)
$ws.Cells.Item(431, 2).Value = 7658985
$ws.Cells.Item(431, 5).Value = "Aguilas Doradas"
$ws.Cells.Item(431, 6).Value = "Fortaleza"
$ws.Cells.Item(431, 7).Value = 1
$ws.Cells.Item(431, 8).Value = 1
$ws.Cells.Item(431, 9).Value = 1
$ws.Cells.Item(431, 10).Value = 0
$ws.Cells.Item(431, 11).Value = "D"
$ws.Cells.Item(431, 12).Value = 1.75
$ws.Cells.Item(431, 13).Value = 3.2
$ws.Cells.Item(431, 14).Value = 5
$ws.Cells.Item(431, 15).Value = 2.05
$ws.Cells.Item(431, 16).Value = 3.2
$ws.Cells.Item(431, 17).Value = 4
$ws.Cells.Item(431, 18).Value = -0.5
$ws.Cells.Item(431, 19).Value = 2.025
$ws.Cells.Item(431, 20).Value = 1.775
$ws.Cells.Item(431, 21).Value = 2
$ws.Cells.Item(431, 22).Value = 1.8
$ws.Cells.Item(431, 23).Value = 2
$ws.Cells.Item(431, 24).Value = -1
$ws.Cells.Item(431, 25).Value = 2.2
$ws.Cells.Item(431, 26).Value = -1
$ws.Cells.Item(431, 27).Value = -1
$ws.Cells.Item(431, 28).Value = 0.7749999999999999
$ws.Cells.Item(431, 29).Value = 0
$ws.Cells.Item(431, 30).Value = 0
